# "Updated Batch addon Step"
# The Msg sheet's "Desc Strt num" test-data row (row 11) is being replaced by a
# new "Desc Strt Spl" scenario: a new shared string is introduced and cell A11
# on the Msg sheet now points at it. The Msg sheet also becomes the active /
# selected sheet (was Login before), with its selection moved to B11.

$wb = $excel.ActiveWorkbook

$msgSheet = $wb.Worksheets.Item("Msg")

# Update the test-data cell: row 11, column A ("Desc Strt num" -> "Desc Strt Spl").
$msgSheet.Range("A11").Value = "Desc Strt Spl"

# Move the active sheet from Login to Msg, and update the in-sheet selection.
$msgSheet.Activate()
$msgSheet.Range("B11").Select()
